# Ajout d'une nouvelle colonne de présence (BR) pour la séance du 28/10/2025
# Mise à jour de l'application

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelle date d'entrainement en tete de colonne BR (ligne 1),
# avec le meme format / alignement que les autres dates (colonne BQ).
$ws.Range("BR1").Value = 45958
$ws.Range("BR1").NumberFormat = "mm-dd-yy"
$ws.Range("BR1").HorizontalAlignment = -4108
$ws.Range("BR1").VerticalAlignment = -4108

# Valeurs de présence pour chaque joueur (ligne 2 a 29), colonne BR
$values = [ordered]@{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "B"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "B"
    19 = "A"
    20 = "P"
    21 = "B"
    22 = "P"
    23 = "RH"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("BR$row")
    $cell.Value = $values[$row]
    $cell.HorizontalAlignment = -4108
}

# La ligne 12 ne comporte pas de seance a cette date (le suivi du joueur
# s'arrete avant la colonne BQ) : on laisse donc BR12 vide.

# Deplace la vue de la feuille pour refleter la nouvelle colonne ajoutee
$ws.Range("BT27").Select()
$ws.Application.ActiveWindow.ScrollColumn = 67
